# lob1004.docx rebuild (2020-10-27 17:38:45 UTC)
#
# The diff shows two things happening near the end of the document:
#   1. A new blank "Normal" paragraph is inserted right before the
#      "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
#   2. The old "(c) 2020 . Contact: ..." paragraph loses its text and is
#      split in two: a plain blank paragraph, followed by a blank
#      paragraph that keeps the original left-justification but now
#      also carries a <w:pageBreakBefore/>.
#
# Everything else in the document is left untouched.

$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs($i)
        if ($candidate.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# --- 1) Insert a new blank paragraph right before "Ver no Jupiter ..." ---
$verIdx = Find-ParagraphContaining $d "Ver no Jupiter Salvar em pdf Salvar em docx"
$verPara = $d.Paragraphs($verIdx)
$verPara.Range.InsertParagraphBefore()

# the freshly inserted paragraph inherited "jc=left" from $verPara; clear it
# so it matches the plain "<w:pStyle .../>"-only blank paragraphs elsewhere.
$blankBeforeVer = $d.Paragraphs($verIdx)
$blankBeforeVer.Alignment = 0

# --- 2) Rework the copyright/footer paragraph ---
$copyIdx = Find-ParagraphContaining $d "Contact: luizeleno@usp.br"
$copyPara = $d.Paragraphs($copyIdx)

# insert a new blank paragraph ahead of it (again strip the inherited jc)
$copyPara.Range.InsertParagraphBefore()
$blankBeforeCopy = $d.Paragraphs($copyIdx)
$blankBeforeCopy.Alignment = 0

# the original paragraph (now shifted one slot later) keeps its jc=left,
# loses its text, and gains a page break before it
$copyPara = $d.Paragraphs($copyIdx + 1)
$copyPara.Range.Text = ""
$copyPara.PageBreakBefore = $true
